$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "Passed"
$ws.Range("H2").Value = "16/04/2021"
$ws.Range("B3").Value = "Yes"
$ws.Range("H3").Value = "16/04/2021"

[void]$ws.Range("B4").Select()
